$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===================== Cell values =====================
$ws.Range("A1").Value = "Название товара`n(обязательное поле)"
$ws.Range("B1").Value = "Код товара`n(до 10 символов)"
$ws.Range("C1").Value = "Артикул `n(до 20 символов)"
$ws.Range("D1").Value = "Штрихкоды  `n(Если оставить штрихкоды пустыми, в терминале будет доступен только визуальный поиск товара)"
$ws.Range("E1").Value = "Отпускная цена (обязательное поле)"
$ws.Range("F1").Value = "Закупочная цена`n(обязательное поле)"
$ws.Range("G1").Value = "Ставка НДС `n(два символа через пробел: `n% (0 – Без НДС)`n“с” - “с рачсетной ставкой”/”н” - “Без”)"
$ws.Range("H1").Value = "Test"
$ws.Range("A2").Value = "Пример"
$ws.Range("C2").Value = "Артикул123456789"
$ws.Range("D2").Value = "Какой-то штрихкод"
$ws.Range("E2").Value = "'18.00"
$ws.Range("F2").Value = "'123.20"
$ws.Range("G2").Value = "0 б"
$ws.Range("B2").Value = 1234567890

# ===================== Row 2 numeric-looking text cleanup (drop quotePrefix style) =====================
$ws.Range("E2").ClearFormats()
$ws.Range("F2").ClearFormats()

# ===================== Header row style (bold Arial 10, center/top, wrap) =====================
$hdr = $ws.Range("A1:H1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 10
$hdr.Font.Bold = $true
$hdr.Font.Color = 0
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.WrapText = $true

# ===================== Rich-text run-level overrides (bold title vs regular note) =====================
$r = $ws.Range("A1").Characters(1,16)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $true
$r = $ws.Range("A1").Characters(17,19)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $false

$r = $ws.Range("B1").Characters(1,11)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $true
$r = $ws.Range("B1").Characters(12,16)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $false

$r = $ws.Range("C1").Characters(1,7)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $true
$r = $ws.Range("C1").Characters(8,18)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $false

$r = $ws.Range("D1").Characters(1,9)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $true
$r = $ws.Range("D1").Characters(10,95)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $false

$r = $ws.Range("E1").Characters(1,14)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $true
$r = $ws.Range("E1").Characters(15,20)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $false

$r = $ws.Range("F1").Characters(1,16)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $true
$r = $ws.Range("F1").Characters(17,19)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $false

$r = $ws.Range("G1").Characters(1,12)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $true
$r = $ws.Range("G1").Characters(13,84)
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Color = 0
$r.Font.Bold = $false

# ===================== I1: styled, empty, center/top (no wrap) =====================
$i1 = $ws.Range("I1")
$i1.HorizontalAlignment = -4108
$i1.VerticalAlignment = -4160

# ===================== Column widths (approximate char widths; engine pixel-snaps) =====================
$ws.Columns.Item(1).ColumnWidth = 26.8333333
$ws.Columns.Item(2).ColumnWidth = 19.5
$ws.Columns.Item(3).ColumnWidth = 22.8333333
$ws.Columns.Item(4).ColumnWidth = 37.5
$ws.Columns.Item(5).ColumnWidth = 16.1666667
$ws.Columns.Item(6).ColumnWidth = 15.1666667
$ws.Columns.Item(7).ColumnWidth = 30.3333333

# ===================== Row heights =====================
$ws.Rows.Item(1).RowHeight = 57

# ===================== Page setup =====================
$ws.PageSetup.Orientation = 1
$ws.PageSetup.HeaderMargin = 36.85
$ws.PageSetup.FooterMargin = 36.85

# ===================== View / selection =====================
$ws.Range("D6").Select()
